$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.673.21'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '2.728.36'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '565.63'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.08'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.69%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.597'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.55%  '
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("E10").Value = '  +4.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.64'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.380'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("D13").Value = '3.210.13'
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.86'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.55%  '
$ws.Range("D15").Value = '63.460.50'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("D17").Value = '2.729.30'
$ws.Range("E17").Value = '  -0.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.59'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.75'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '354.51'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.57'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.26%  '
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.522'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.57'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.87%  '
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.38'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").Value = '0.0₃0908'
$ws.Range("E28").Value = '  +0.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.98'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.18'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.33'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +11.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '166.22'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.92'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.06'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.34%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.48'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.18%  '
$ws.Range("B36").Value = 'USDe'
$ws.Range("C36").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.81'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.975'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '347.50'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +6.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.30'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.16%  '
$ws.Range("E41").Value = '  -0.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.68'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.92'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.04'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0584'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.21%  '
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0251'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.100'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.47%  '
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.84'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.46%  '
$ws.Range("E51").Value = '  +0.12%  '
